$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the "Vega Monumental Concepción" weekly block
# (new price report week, dated 44461), pushing the existing rows 319-346 down
# to 322-349.
$ws.Range("A319:R321").EntireRow.Insert()

# Row 319: Conconina(o)
$ws.Range("A319").Value = 11
$ws.Range("B319").Value = "Vega Monumental Concepción"
$ws.Range("C319").Value = "Bíobío"
$ws.Range("D319").Value = 44461
$ws.Range("E319").Value = 8
$ws.Range("F319").Value = 100112033
$ws.Range("G319").Value = "Lechuga"
$ws.Range("H319").Value = "Conconina(o)"
$ws.Range("I319").Value = "Primera"
$ws.Range("J319").Value = 100
$ws.Range("K319").Value = 5500
$ws.Range("L319").Value = 6000
$ws.Range("M319").Value = 5750
$ws.Range("N319").Value = "`$/caja 10 unidades"
$ws.Range("O319").Value = "Región Metropolitana"
$ws.Range("P319").Value = 575
$ws.Range("Q319").Value = 10
$ws.Range("R319").Value = "Hortaliza"

# Row 320: Francesa morada
$ws.Range("A320").Value = 11
$ws.Range("B320").Value = "Vega Monumental Concepción"
$ws.Range("C320").Value = "Bíobío"
$ws.Range("D320").Value = 44461
$ws.Range("E320").Value = 8
$ws.Range("F320").Value = 100112033
$ws.Range("G320").Value = "Lechuga"
$ws.Range("H320").Value = "Francesa morada"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 100
$ws.Range("K320").Value = 6000
$ws.Range("L320").Value = 6500
$ws.Range("M320").Value = 6250
$ws.Range("N320").Value = "`$/caja 15 unidades"
$ws.Range("O320").Value = "Región Metropolitana"
$ws.Range("P320").Value = 417
$ws.Range("Q320").Value = 15
$ws.Range("R320").Value = "Hortaliza"

# Row 321: Marina
$ws.Range("A321").Value = 11
$ws.Range("B321").Value = "Vega Monumental Concepción"
$ws.Range("C321").Value = "Bíobío"
$ws.Range("D321").Value = 44461
$ws.Range("E321").Value = 8
$ws.Range("F321").Value = 100112033
$ws.Range("G321").Value = "Lechuga"
$ws.Range("H321").Value = "Marina"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 100
$ws.Range("K321").Value = 5500
$ws.Range("L321").Value = 6000
$ws.Range("M321").Value = 5750
$ws.Range("N321").Value = "`$/caja 15 unidades"
$ws.Range("O321").Value = "Región Metropolitana"
$ws.Range("P321").Value = 383
$ws.Range("Q321").Value = 15
$ws.Range("R321").Value = "Hortaliza"
